$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.505.47"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.573.65"
$ws.Range("E3").Value = "  -2.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "2.583.23"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0996"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.329"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.133"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "3.025.10"
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "58.484.97"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "2.566.92"
$ws.Range("E17").Value = "  -4.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000131"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "335.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.401"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "0.0₃0700"
$ws.Range("E30").Value = "  -11.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.64%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.822"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.995"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "271.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.584"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0943"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0514"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("D48").Value = "1.963.57"
$ws.Range("E48").Value = "  -3.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0218"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.54%  "
